$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '29.837.11'
$ws.Range("E2").Value = '  -0.24%  '
Set-TextValue "D3" '1.876.13'
$ws.Range("E3").Value = '  +0.04%  '
Set-TextValue "D4" '1.003'
$ws.Range("E4").Value = '  +0.32%  '
Set-TextValue "D5" '0.7200'
$ws.Range("E5").Value = '  -2.64%  '
Set-TextValue "D6" '242.26'
$ws.Range("E7").Value = '  +0.26%  '
Set-TextValue "D8" '0.3148'
$ws.Range("E8").Value = '  -0.51%  '
Set-TextValue "D9" '0.07346'
$ws.Range("E9").Value = '  +1.84%  '
Set-TextValue "D10" '24.55'
$ws.Range("E10").Value = '  -0.92%  '
Set-TextValue "D11" '0.08195'
$ws.Range("E11").Value = '  -2.44%  '
Set-TextValue "D12" '0.7446'
$ws.Range("E12").Value = '  -0.80%  '
Set-TextValue "D13" '1.880.60'
$ws.Range("E13").Value = '  -0.35%  '
Set-TextValue "D14" '5.329'
$ws.Range("E14").Value = '  -1.82%  '
Set-TextValue "D15" '92.52'
$ws.Range("E15").Value = '  -0.06%  '
Set-TextValue "D16" '29.882.63'
$ws.Range("E16").Value = '  -0.07%  '
Set-TextValue "D17" '6.015'
$ws.Range("E17").Value = '  -1.31%  '
Set-TextValue "D18" '246.76'
$ws.Range("E18").Value = '  +1.29%  '
Set-TextValue "D19" '13.50'
$ws.Range("E19").Value = '  -0.51%  '
Set-TextValue "D20" '0.000007891'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("E21").Value = '  +0.21%  '
Set-TextValue "D22" '2.139.96'
$ws.Range("E22").Value = '  +0.79%  '
Set-TextValue "D23" '1.003'
$ws.Range("E23").Value = '  -0.16%  '
Set-TextValue "D24" '7.723'
$ws.Range("E24").Value = '  -3.09%  '
Set-TextValue "D25" '9.246'
$ws.Range("E25").Value = '  -0.35%  '
Set-TextValue "D26" '0.1500'
$ws.Range("E26").Value = '  -3.31%  '
Set-TextValue "D27" '164.16'
$ws.Range("E28").Value = '  -0.13%  '
Set-TextValue "D29" '2.008'
$ws.Range("E29").Value = '  -1.24%  '
Set-TextValue "D30" '1.426'
$ws.Range("E30").Value = '  -5.48%  '
Set-TextValue "D31" '4.533'
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("E32").Value = '  -0.27%  '
Set-TextValue "D33" '4.173'
$ws.Range("E33").Value = '  -2.26%  '
Set-TextValue "D34" '0.05452'
$ws.Range("E34").Value = '  +2.61%  '
Set-TextValue "D35" '1.229'
$ws.Range("E35").Value = '  -0.67%  '
Set-TextValue "D36" '0.7341'
$ws.Range("E36").Value = '  -2.69%  '
Set-TextValue "D37" '1.002'
$ws.Range("E37").Value = '  +0.46%  '
Set-TextValue "D38" '2.701'
$ws.Range("E38").Value = '  +0.02%  '
Set-TextValue "D39" '0.01914'
$ws.Range("E39").Value = '  -2.36%  '
Set-TextValue "D40" '2.744'
$ws.Range("E40").Value = '  -0.36%  '
Set-TextValue "D41" '0.4445'
$ws.Range("E41").Value = '  -1.97%  '
Set-TextValue "D42" '0.8971'
$ws.Range("E42").Value = '  +4.68%  '
Set-TextValue "D43" '6.007'
$ws.Range("E43").Value = '  -0.45%  '
Set-TextValue "D44" '71.54'
$ws.Range("E44").Value = '  -1.35%  '

# Row 45/46: Maker and PaxDollar swap positions with updated values
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D45" '1.003'
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D46" '1.040.02'
$ws.Range("E46").Value = '  -6.91%  '

Set-TextValue "D47" '103.67'
$ws.Range("E47").Value = '  +0.34%  '
Set-TextValue "D48" '7.461'
$ws.Range("E48").Value = '  -2.27%  '
Set-TextValue "D49" '1.807'
$ws.Range("E49").Value = '  -1.79%  '
Set-TextValue "D50" '9.628'
$ws.Range("E50").Value = '  +0.96%  '
Set-TextValue "D51" '2.027.60'
$ws.Range("E51").Value = '  +0.33%  '
